# Split the "OGGETTO: " run into a bold "OGGETTO:" run followed by a
# (still italic) run containing just the trailing space, matching the
# author's edit that bolded the "OGGETTO:" label.

$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "OGGETTO: "
$found = $rng.Find.Execute($null, $false, $false, $false, $false, $false, `
                            $true, 1, $false, $null, 0)

if ($found) {
    $start = $rng.Start
    $end = $rng.End

    # Range covering just "OGGETTO:" (without the trailing space).
    $boldRange = $d.Range($start, $end - 1)

    $boldRange.Font.Bold = $true
    $boldRange.Font.Italic = $false
}
